$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 143, shifting existing rows 143:214 down to 144:215.
$ws.Rows.Item(143).Insert()

# Populate the newly inserted row 143 with the new weekly record.
$ws.Cells.Item(143, 1).Value = 8
$ws.Cells.Item(143, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(143, 3).Value = "Coquimbo"
$ws.Cells.Item(143, 4).Value = 44523
$ws.Cells.Item(143, 5).Value = 4
$ws.Cells.Item(143, 6).Value = 100112032
$ws.Cells.Item(143, 7).Value = "Zapallo italiano"
$ws.Cells.Item(143, 8).Value = "Sin especificar"
$ws.Cells.Item(143, 9).Value = "Primera"
$ws.Cells.Item(143, 10).Value = 500
$ws.Cells.Item(143, 11).Value = 10000
$ws.Cells.Item(143, 12).Value = 11000
$ws.Cells.Item(143, 13).Value = 10500
$ws.Cells.Item(143, 14).Value = "$/caja 70 unidades"
$ws.Cells.Item(143, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(143, 16).Value = 150
$ws.Cells.Item(143, 17).Value = 70
$ws.Cells.Item(143, 18).Value = "Hortaliza"
